$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 26 ("Added menu screen with info") above the existing
# "Cleanup code" row, reusing the formatting of a typical DONE task row
# (row 2) so the new row gets the same styles (time format on B, status
# fill on C) as the rest of the table.
$ws.Rows("2:2").Copy()
$ws.Rows("26:26").Insert()

# Fill in the new row's own content.
$ws.Range("A26").Value = "Added menu screen with info"
$ws.Range("B26").Value = 0.013888888888888888
$ws.Range("C26").Value = "DONE"

# Restore the selection, now pointing at A28 instead of C28.
$null = $ws.Range("A28").Select()
